# Trade #5 closed at 2026-02-17 12:26:56 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.99   # Current Capital
$summary.Range("B4").Value = -0.01     # Total P&L $
$summary.Range("B5").Value = -0.04     # Total P&L %
$summary.Range("B6").Value = 5         # Total Trades
$summary.Range("B8").Value = 3         # Losing Trades
$summary.Range("B9").Value = 40        # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.99      # Capital (MarketMaking row)
$status.Range("D4").Value = 5          # Trades
$status.Range("E4").Value = -0.01      # P&L $
$status.Range("F4").Value = -0.01      # P&L %
$status.Range("G4").Value = 40         # Win Rate %

# ---------------------------------------------------------------------------
# New trade row (Trade #5), appended to both the "All Trades" sheet and the
# per-strategy "MarketMaking" sheet.
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A6").Value = 5
    # Force column B to be stored as plain text instead of being
    # auto-converted to a date serial number.
    $ws.Range("B6").NumberFormat = "@"
    $ws.Range("B6").Value = "2026-02-17"
    $ws.Range("C6").Value = "12:26:50"
    $ws.Range("D6").Value = "MarketMaking"
    $ws.Range("E6").Value = "UP"
    $ws.Range("F6").Value = 0.16
    $ws.Range("G6").Value = 0.14313
    $ws.Range("H6").Value = "CLOSED"
    $ws.Range("I6").Value = -10.5439
    $ws.Range("J6").Value = -0.02
    $ws.Range("K6").Value = 99.99
    $ws.Range("L6").Value = 0
    $ws.Range("M6").Value = 0
    $ws.Range("N6").Value = 0.6
    $ws.Range("O6").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P6").Value = "early_exit"
    $ws.Range("Q6").Value = 0.1
}
